$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 13, shifting existing rows 13:100 down to 14:101
$ws.Rows("13:13").Insert()

# Populate the newly inserted row 13 with the new record's data
$ws.Range("A13").Value = 10
$ws.Range("B13").Value = 'Vega Modelo de Temuco'
$ws.Range("C13").Value = 'La Araucanía'
$ws.Range("D13").Value = 45022
$ws.Range("E13").Value = 9
$ws.Range("F13").Value = 'Fruta'
$ws.Range("G13").Value = 100108
$ws.Range("H13").Value = 'Tropicales y subtropicales'
$ws.Range("I13").Value = 100108004
$ws.Range("J13").Value = 'Papaya'
$ws.Range("K13").Value = 'Cultivar IV Región'
$ws.Range("L13").Value = 'Primera'
$ws.Range("M13").Value = 150
$ws.Range("N13").Value = 27000
$ws.Range("O13").Value = 27000
$ws.Range("P13").Value = 27000
$ws.Range("Q13").Value = '$/bandeja 10 kilos'
$ws.Range("R13").Value = 'Provincia del Elquí'
$ws.Range("S13").Value = 2700
$ws.Range("T13").Value = 10
